# Auto-generated edit script applying weekly price/date updates
# for rows 13-35 (Hortaliza, Macroferia Regional de Talca - Esparragos)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 13
$ws.Range("D13").Value = 44473
$ws.Range("J13").Value = 4000
$ws.Range("K13").Value = 1200
$ws.Range("L13").Value = 1200
$ws.Range("M13").Value = 1200
$ws.Range("O13").Value = 'Provincia de Linares'
$ws.Range("P13").Value = 1200

# Row 14
$ws.Range("D14").Value = 44512
$ws.Range("J14").Value = 5000
$ws.Range("K14").Value = 800
$ws.Range("L14").Value = 800
$ws.Range("M14").Value = 800
$ws.Range("O14").Value = 'Región del Maule'
$ws.Range("P14").Value = 800

# Row 15
$ws.Range("D15").Value = 44460
$ws.Range("J15").Value = 2000
$ws.Range("K15").Value = 2000
$ws.Range("L15").Value = 2000
$ws.Range("M15").Value = 2000
$ws.Range("O15").Value = 'Provincia de Linares'
$ws.Range("P15").Value = 2000

# Row 16
$ws.Range("D16").Value = 44489
$ws.Range("J16").Value = 4000
$ws.Range("K16").Value = 900
$ws.Range("L16").Value = 900
$ws.Range("M16").Value = 900
$ws.Range("P16").Value = 900

# Row 17
$ws.Range("D17").Value = 44497
$ws.Range("J17").Value = 5000

# Row 18
$ws.Range("D18").Value = 44482
$ws.Range("J18").Value = 4000
$ws.Range("L18").Value = 1000
$ws.Range("M18").Value = 950
$ws.Range("P18").Value = 950

# Row 19
$ws.Range("D19").Value = 44516
$ws.Range("J19").Value = 3000
$ws.Range("K19").Value = 1000
$ws.Range("M19").Value = 1000
$ws.Range("O19").Value = 'Provincia de Limarí'
$ws.Range("P19").Value = 1000

# Row 20
$ws.Range("D20").Value = 44175
$ws.Range("J20").Value = 800
$ws.Range("L20").Value = 1100
$ws.Range("M20").Value = 1050
$ws.Range("O20").Value = 'Provincia de Linares'
$ws.Range("P20").Value = 1050

# Row 21
$ws.Range("D21").Value = 44168
$ws.Range("J21").Value = 3000
$ws.Range("L21").Value = 1000
$ws.Range("M21").Value = 1000
$ws.Range("P21").Value = 1000

# Row 22
$ws.Range("D22").Value = 44475
$ws.Range("J22").Value = 5000
$ws.Range("L22").Value = 1100
$ws.Range("M22").Value = 1040
$ws.Range("P22").Value = 1040

# Row 23
$ws.Range("D23").Value = 44483
$ws.Range("J23").Value = 4000
$ws.Range("K23").Value = 900
$ws.Range("L23").Value = 1000
$ws.Range("M23").Value = 950
$ws.Range("P23").Value = 950

# Row 24
$ws.Range("D24").Value = 44162
$ws.Range("K24").Value = 1000
$ws.Range("M24").Value = 1000
$ws.Range("N24").Value = '$/atado'
$ws.Range("P24").Value = 1000

# Row 25
$ws.Range("D25").Value = 44467
$ws.Range("J25").Value = 2000
$ws.Range("K25").Value = 1800
$ws.Range("L25").Value = 1800
$ws.Range("M25").Value = 1800
$ws.Range("N25").Value = '$/kilo'
$ws.Range("P25").Value = 1800

# Row 26
$ws.Range("D26").Value = 44167
$ws.Range("K26").Value = 1000
$ws.Range("L26").Value = 1000
$ws.Range("M26").Value = 1000
$ws.Range("O26").Value = 'Región del Maule'
$ws.Range("P26").Value = 1000

# Row 27
$ws.Range("D27").Value = 44496
$ws.Range("J27").Value = 4000
$ws.Range("K27").Value = 900
$ws.Range("L27").Value = 900
$ws.Range("M27").Value = 900
$ws.Range("O27").Value = 'Provincia de Linares'
$ws.Range("P27").Value = 900

# Row 28
$ws.Range("D28").Value = 44498
$ws.Range("J28").Value = 5000

# Row 29
$ws.Range("D29").Value = 44477
$ws.Range("J29").Value = 4000
$ws.Range("K29").Value = 1000
$ws.Range("L29").Value = 1000
$ws.Range("M29").Value = 1000
$ws.Range("P29").Value = 1000

# Row 30
$ws.Range("D30").Value = 44487
$ws.Range("J30").Value = 5000
$ws.Range("K30").Value = 800
$ws.Range("L30").Value = 800
$ws.Range("M30").Value = 800
$ws.Range("P30").Value = 800

# Row 31
$ws.Range("D31").Value = 44505
$ws.Range("J31").Value = 6000

# Row 32
$ws.Range("D32").Value = 44509

# Row 33
$ws.Range("D33").Value = 44494
$ws.Range("J33").Value = 4000
$ws.Range("K33").Value = 900
$ws.Range("L33").Value = 900
$ws.Range("M33").Value = 900
$ws.Range("P33").Value = 900

# Row 34
$ws.Range("D34").Value = 44469
$ws.Range("J34").Value = 3000
$ws.Range("K34").Value = 1200
$ws.Range("L34").Value = 1200
$ws.Range("M34").Value = 1200
$ws.Range("P34").Value = 1200

# Row 35
$ws.Range("D35").Value = 44518
$ws.Range("K35").Value = 1000
$ws.Range("L35").Value = 1000
$ws.Range("M35").Value = 1000
$ws.Range("P35").Value = 1000
